# Remove "laudos" rows from email export sheet; rebuild header + data per new schema
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force a numeric-looking string to be stored as text (shared string),
    # avoiding NumberFormat changes that would add new style records.
    $helper = $ws.Range("ZZ1")
    $helper.Formula = "=" + [char]34 + $text + [char]34
    $helper.Copy()
    $range.PasteSpecial(-4163)
    $helper.Clear()
}

# Header row (row 1)
$ws.Range("B1").Value = "Nome Cliente"
$ws.Range("C1").Value = "Endereço"
$ws.Range("D1").Value = "Numero"
$ws.Range("E1").Value = "Complemento"
$ws.Range("F1").Value = "Bairro"
$ws.Range("G1").Value = "CEP"
$ws.Range("H1").Value = "Município"
$ws.Range("I1").Value = "UF"
$ws.Range("J1").Value = "Tipo"
$ws.Range("K1").Value = "Lead ID"
$ws.Range("L1").Value = "Observações"
$ws.Range("M1").Value = "Endereço Completo"

# M1 needs the same bold/border header style as the rest of row 1 (xlPasteFormats)
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Leticia Maria Pereira Rodrigues"
$ws.Range("C2").Value = "Alameda Dos Guainumbis"
Set-TextValue $ws.Range("D2") "261"
$ws.Range("F2").Value = "Alameda dos Guainumbis - Planalto Paulista"
Set-TextValue $ws.Range("G2") "4067000"
$ws.Range("H2").Value = "São Paulo"
$ws.Range("I2").Value = "SP"
$ws.Range("J2").Value = "Casa Residencial"
$ws.Range("K2").Value = "44YX64N"
$ws.Range("L2").Value = "Não foi encontrado nada sobre esse bairro no Wikpedia"
$ws.Range("M2").Value = "Alameda Dos Guainumbis 261-"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Leticia De Lima Purificação"
$ws.Range("C3").Value = "Rua Vitório Favalli"
Set-TextValue $ws.Range("D3") "57"
$ws.Range("F3").Value = "Vila Maria de Maggi"
Set-TextValue $ws.Range("G3") "8680120"
$ws.Range("H3").Value = "Suzano"
$ws.Range("I3").Value = "SP"
$ws.Range("J3").Value = "Casa Residencial"
$ws.Range("K3").Value = "2A4JF4Q"
$ws.Range("L3").Value = "Não foi encontrado nada sobre esse bairro no Wikpedia"
$ws.Range("M3").Value = "Rua Vitório Favalli 57-"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Gabriel Felipe Guimarães Coutinho Cortez"
$ws.Range("C4").Value = "Sds Bloco F E G"
$ws.Range("D4").Value = "Lote 41 S/n"
$ws.Range("E4").Value = "edificio Eldorado Sala 609"
$ws.Range("F4").Value = "Asa Sul"
Set-TextValue $ws.Range("G4") "70392900"
$ws.Range("H4").Value = "Brasilia"
$ws.Range("I4").Value = "DF"
$ws.Range("J4").Value = "Sala Comercial"
$ws.Range("K4").Value = "HNXQSGX"
$ws.Range("L4").Value = "Brasília (AFI: /bɾaˈziljɐ/ ou AFI: /bɾaˈziʎɐ/) é a capital federal do Brasil e a sede de governo do Distrito Federal. A capital está localizada na região Centro-Oeste do país, ao longo da região geográfica conhecida como Planalto Central. Segundo estimativa do Instituto Brasileiro de Geografia e Estatística (IBGE) para 2018, sua população era de 2 974 703 habitantes (4 284 676 em sua área metropolitana), sendo, então, a terceira cidade mais populosa do país. Brasília é também a quinta concentração urbana mais populosa do Brasil. A capital brasileira é a maior cidade do mundo construída no século XX.A cidade possui o maior produto interno bruto per capita em relação às capitais, o quarto maior entre as principais cidades da América Latina e cerca de três vezes maior que a renda média brasileira. Como capital nacional, Brasília abriga a sede dos três poderes da República (Executivo, Legislativo e Judiciário) e 127 embaixadas estrangeiras. A política de planejamento da cidade, como a localização de prédios residenciais em grandes áreas urbanas, a construção da cidade através de enormes avenidas e a sua divisão em setores, tem provocado debates sobre o estilo de vida nas grandes cidades no século XX. O projeto da cidade a divide em blocos numerados, além de setores para atividades pré-determinadas, como o Setor Hoteleiro, Bancário ou de Embaixadas."
$ws.Range("M4").Value = "Sds Bloco F E G Lote Lote 41 S/n-edificio Eldorado Sala 609"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Sandra Ferreira Portugal Moura"
$ws.Range("C5").Value = "Rua Guilherme De Almeida"
Set-TextValue $ws.Range("D5") "456"
$ws.Range("E5").Value = "apto 602"
$ws.Range("F5").Value = "Santo Antônio"
Set-TextValue $ws.Range("G5") "30350230"
$ws.Range("H5").Value = "Belo Horizonte"
$ws.Range("I5").Value = "MG"
$ws.Range("J5").Value = "Apartamento"
$ws.Range("K5").Value = "HEGLGUH"
$ws.Range("L5").Value = "Esta é uma lista contendo os todos os bairros e territórios do município de Belo Horizonte, capital do estado de Minas Gerais. Conforme levantamento da Embrapa no `"trabalho de identificação, mapeamento e quantificação das áreas urbanas do Brasil`", em 2015 Belo Horizonte possuía area urbana de 314 km². Assim, considerando a área total do município, que é de 331 km², em relação a sua área urbana, o município tem 95% de seu território urbanizado.Segundo informações da a prefeitura da cidade, em 2021 Belo Horizonte possui ao todo 487 bairros, onde estão distribuidos 15.992 logradouros, dentro os quais há exatamente 11.479 ruas.  Os demais tipos de logradoruros incluem avenidas, estradas, túneis, becos, dentre outros, sendo que destes, os becos totalizam 2.635, número em constante constante alteração a medida em que novas ocupações vão sendo descobertas no município."
$ws.Range("M5").Value = "Rua Guilherme De Almeida 456-apto 602"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Leticia De Lima Purificação"
$ws.Range("C6").Value = "Rua Vitório Favalli"
Set-TextValue $ws.Range("D6") "57"
$ws.Range("F6").Value = "Vila Maria de Maggi"
Set-TextValue $ws.Range("G6") "8680120"
$ws.Range("H6").Value = "Suzano"
$ws.Range("I6").Value = "SP"
$ws.Range("J6").Value = "Casa Residencial"
$ws.Range("K6").Value = "2A4JF4Q"
$ws.Range("L6").Value = "Não foi encontrado nada sobre esse bairro no Wikpedia"
$ws.Range("M6").Value = "Rua Vitório Favalli 57-"

# A2:A6 all carry the bold/border "id column" style (xlPasteFormats from A2)
$ws.Range("A2").Copy()
$ws.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
